$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 297
$ws.Range("I6").Value = 251
$ws.Range("J6").Value = 320
$ws.Range("K6").Value = 753
$ws.Range("L6").Value = 960
$ws.Range("M6").Value = -641
$ws.Range("N6").Value = -1184
$ws.Range("H11").Value = 610.2105
$ws.Range("I11").Value = 610.2105
$ws.Range("K11").Value = 610.2105
$ws.Range("M11").Value = -470.2105
$ws.Range("H33").Value = 370.95
$ws.Range("I33").Value = 295
$ws.Range("K33").Value = 295
$ws.Range("M33").Value = -66
$ws.Range("H116").Value = 9218.272000000001
$ws.Range("I116").Value = 7999.1665
$ws.Range("J116").Value = 10681.2
$ws.Range("K116").Value = 7999.1665
$ws.Range("L116").Value = 10681.2
$ws.Range("M116").Value = -4557.1665
$ws.Range("N116").Value = -17565.2
$ws.Range("H127").Value = 2461.2
$ws.Range("J127").Value = 3949.8
$ws.Range("L127").Value = 11849.4
$ws.Range("N127").Value = -21769.4
$ws.Range("H129").Value = 1776.3
$ws.Range("J129").Value = 2899.4
$ws.Range("L129").Value = 8698.200000000001
$ws.Range("N129").Value = -18698.2
$ws.Range("H135").Value = 3104.3635
$ws.Range("I135").Value = 2693.3928
$ws.Range("J135").Value = 5405.8
$ws.Range("K135").Value = 24240.5352
$ws.Range("L135").Value = 48652.2
$ws.Range("M135").Value = -21705.5352
$ws.Range("N135").Value = -53722.2
$ws.Range("H138").Value = 6150.8184
$ws.Range("J138").Value = 6563.1895
$ws.Range("L138").Value = 19689.5685
$ws.Range("N138").Value = -29969.5685

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12250.298
$ws.Range("I32").Value = 9766.75
$ws.Range("J32").Value = 48675.668
$ws.Range("K32").Value = 9766.75
$ws.Range("L32").Value = 48675.668
$ws.Range("M32").Value = -9479.75
$ws.Range("N32").Value = -49249.668
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H45").Value = 1708.875
$ws.Range("I45").Value = 1134.2307
$ws.Range("K45").Value = 1134.2307
$ws.Range("M45").Value = -757.2307000000001
$ws.Range("H101").Value = 44975
$ws.Range("J101").Value = 44975
$ws.Range("L101").Value = 44975
$ws.Range("N101").Value = -51465
$ws.Range("H122").Value = 4844.778
$ws.Range("I122").Value = 4616.625
$ws.Range("J122").Value = 5027.3
$ws.Range("K122").Value = 13849.875
$ws.Range("L122").Value = 15081.9
$ws.Range("M122").Value = -11399.875
$ws.Range("N122").Value = -19981.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H128").Value = 2600
$ws.Range("I128").Value = 2600
$ws.Range("K128").Value = 7800
$ws.Range("M128").Value = -5310
$ws.Range("H135").Value = 57624.715
$ws.Range("J135").Value = 61395.5
$ws.Range("L135").Value = 61395.5
$ws.Range("N135").Value = -71535.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27737.773
$ws.Range("I31").Value = 2942.3333
$ws.Range("J31").Value = 57492.3
$ws.Range("K31").Value = 2942.3333
$ws.Range("L31").Value = 57492.3
$ws.Range("M31").Value = -2647.3333
$ws.Range("N31").Value = -58082.3
$ws.Range("H34").Value = 27737.773
$ws.Range("I34").Value = 2942.3333
$ws.Range("J34").Value = 57492.3
$ws.Range("K34").Value = 2942.3333
$ws.Range("L34").Value = 57492.3
$ws.Range("M34").Value = -2740.3333
$ws.Range("N34").Value = -57896.3
$ws.Range("H104").Value = 52529.5
$ws.Range("I104").Value = 52529.5
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 52529.5
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -49908.5
$ws.Range("N104").ClearContents()
$ws.Range("H105").Value = 4469.857
$ws.Range("I105").Value = 1749
$ws.Range("K105").Value = 1749
$ws.Range("M105").Value = -2
$ws.Range("H108").Value = 82484.664
$ws.Range("J108").Value = 82484.664
$ws.Range("L108").Value = 82484.664
$ws.Range("N108").Value = -90164.664
$ws.Range("H131").Value = 34065.332
$ws.Range("I131").Value = 27598
$ws.Range("J131").Value = 47000
$ws.Range("K131").Value = 27598
$ws.Range("L131").Value = 47000
$ws.Range("M131").Value = -22558
$ws.Range("N131").Value = -57080
$ws.Range("H134").Value = 2032.7028
$ws.Range("I134").Value = 1594.2333
$ws.Range("K134").Value = 4782.699900000001
$ws.Range("M134").Value = -2247.699900000001
$ws.Range("H141").Value = 549903.2
$ws.Range("J141").Value = 585610.5600000001
$ws.Range("L141").Value = 585610.5600000001
$ws.Range("N141").Value = -595970.5600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 25100.666
$ws.Range("J2").Value = 31647.21
$ws.Range("L2").Value = 189883.26
$ws.Range("N2").Value = -190109.26
$ws.Range("H5").Value = 1334255.2
$ws.Range("I5").Value = 784.375
$ws.Range("K5").Value = 2353.125
$ws.Range("M5").Value = -2241.125
$ws.Range("H51").Value = 50
$ws.Range("I51").Value = 50
$ws.Range("K51").Value = 150
$ws.Range("M51").Value = 310
$ws.Range("H80").Value = 10500
$ws.Range("I80").Value = 9998.5
$ws.Range("J80").Value = 11001.5
$ws.Range("K80").Value = 29995.5
$ws.Range("L80").Value = 33004.5
$ws.Range("M80").Value = -29059.5
$ws.Range("N80").Value = -34876.5
$ws.Range("H83").Value = 10500
$ws.Range("I83").Value = 9998.5
$ws.Range("J83").Value = 11001.5
$ws.Range("K83").Value = 89986.5
$ws.Range("L83").Value = 99013.5
$ws.Range("M83").Value = -85306.5
$ws.Range("N83").Value = -108373.5
$ws.Range("H113").Value = 3128.3333
$ws.Range("J113").Value = 3208.182
$ws.Range("L113").Value = 9624.545999999998
$ws.Range("N113").Value = -13964.546
$ws.Range("H122").Value = 3124.5454
$ws.Range("J122").Value = 4302.857
$ws.Range("L122").Value = 38725.713
$ws.Range("N122").Value = -43625.713
$ws.Range("H135").Value = 1334255.2
$ws.Range("I135").Value = 784.375
$ws.Range("K135").Value = 7059.375
$ws.Range("M135").Value = -4524.375
$ws.Range("H136").Value = 3116.28
$ws.Range("I136").Value = 3041.9583
$ws.Range("J136").Value = 4900
$ws.Range("K136").Value = 9125.874899999999
$ws.Range("L136").Value = 14700
$ws.Range("M136").Value = -4025.874899999999
$ws.Range("N136").Value = -24900
$ws.Range("H137").Value = 75189.5
$ws.Range("I137").Value = 2415.2
$ws.Range("J137").Value = 115619.664
$ws.Range("K137").Value = 7245.599999999999
$ws.Range("L137").Value = 346858.992
$ws.Range("M137").Value = -2145.599999999999
$ws.Range("N137").Value = -357058.992
$ws.Range("H139").Value = 4486.467
$ws.Range("I139").Value = 3223.8823
$ws.Range("J139").Value = 6137.5386
$ws.Range("K139").Value = 9671.6469
$ws.Range("L139").Value = 18412.6158
$ws.Range("M139").Value = -4531.6469
$ws.Range("N139").Value = -28692.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5000000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 5000000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 1697.5227
$ws.Range("I97").Value = 817.4074000000001
$ws.Range("K97").Value = 817.4074000000001
$ws.Range("M97").Value = -321.4074000000001
$ws.Range("H122").Value = 3674.0645
$ws.Range("I122").Value = 3183.7778
$ws.Range("J122").Value = 6983.5
$ws.Range("K122").Value = 9551.3334
$ws.Range("L122").Value = 20950.5
$ws.Range("M122").Value = -7101.3334
$ws.Range("N122").Value = -25850.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 2977.8572
$ws.Range("I30").Value = 2977.8572
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2977.8572
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2869.8572
$ws.Range("N30").ClearContents()
$ws.Range("H82").Value = 2806.5833
$ws.Range("J82").Value = 3751
$ws.Range("L82").Value = 3751
$ws.Range("N82").Value = -4473
$ws.Range("H85").Value = 2806.5833
$ws.Range("J85").Value = 3751
$ws.Range("L85").Value = 3751
$ws.Range("N85").Value = -6247
$ws.Range("H131").Value = 200000
$ws.Range("J131").Value = 200000
$ws.Range("L131").Value = 200000
$ws.Range("N131").Value = -210080
$ws.Range("H132").Value = 7012.8
$ws.Range("I132").Value = 5654.778
$ws.Range("J132").Value = 9049.833000000001
$ws.Range("K132").Value = 16964.334
$ws.Range("L132").Value = 27149.499
$ws.Range("M132").Value = -14434.334
$ws.Range("N132").Value = -32209.499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11644
$ws.Range("I81").Value = 10525
$ws.Range("J81").Value = 15001
$ws.Range("K81").Value = 21050
$ws.Range("L81").Value = 30002
$ws.Range("M81").Value = -19989
$ws.Range("N81").Value = -32124
$ws.Range("H84").Value = 11644
$ws.Range("I84").Value = 10525
$ws.Range("J84").Value = 15001
$ws.Range("K84").Value = 105250
$ws.Range("L84").Value = 150010
$ws.Range("M84").Value = -99946
$ws.Range("N84").Value = -160618
$ws.Range("H136").Value = 6182
$ws.Range("I136").Value = 3066
$ws.Range("K136").Value = 9198
$ws.Range("M136").Value = -6648
